$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New summary formulas comparing pre/post groups (added in column S/T,
# rows 21-26) on the "p3_sec_outcomes_kyle" sheet.
$ws.Range("S21").Formula = "=AVERAGE(C2:C5)+AVERAGE(C15:C25)"
$ws.Range("S22").Formula = "=AVERAGE(C6:C14)+AVERAGE(C26:C34)"
$ws.Range("S23").Formula = "=S22-S21"
$ws.Range("T23").Formula = "=_xlfn.STDEV.P(C2:C5, C15:C25)"
$ws.Range("S24").Formula = "=S23/T26"
$ws.Range("T24").Formula = "=_xlfn.STDEV.P(C6:C14, C26:C34)"
$ws.Range("T25").Formula = "=(POWER(T23, 2) + POWER(T24, 2))/2"
$ws.Range("T26").Formula = "=SQRT(T25)"

# Selection left on T29 when the file was last saved.
$ws.Range("T29").Select()
